$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph (plain
# substring search -- avoids regex-escaping the parentheses/accents).
$target = $null
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Docente(s)")) {
        $target = $p
        $targetIndex = $i
        break
    }
    $i++
}

if ($target -eq $null) {
    throw "Could not locate 'Docente(s) Responsável(eis)' paragraph"
}

# Insert a fresh (empty) paragraph right after the heading paragraph.
$insertionPoint = $target.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

# Grab the whole new (still empty) paragraph -- including its end-of-paragraph
# mark -- and replace it in one shot with the real OOXML for the bullet list
# of professors, so each "name + line break" combo keeps its own <w:r> run,
# matching how the rest of the document (e.g. "Créditos-aula") is authored.
$newPara = $d.Paragraphs($targetIndex + 1)
$fullRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
       '<w:r><w:t>471420 - Carlos Antonio Reis Pereira Baptista</w:t><w:br/></w:r>' +
       '<w:r><w:t>3480026 - João Paulo Pascon</w:t><w:br/></w:r>' +
       '<w:r><w:t>5840793 - Sérgio Schneider</w:t><w:br/></w:r>' +
       '<w:r><w:t>7797767 - Viktor Pastoukhov</w:t></w:r>' +
       '</w:p>'

[void]$fullRange.InsertXML($xml)

Write-Output "Inserted professors bullet list after the Docente(s) Responsavel(eis) heading paragraph."
